# Regenerate save_data: replace column G (Strike# -> K) values for rows 2-13
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 2
    6  = 3
    7  = 3
    8  = 5
    9  = 4
    10 = 4
    11 = 5
    12 = 6
    13 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
